# Rebuild the document body to match the target revision.
$apos = [char]0x2019

$body = "<w:p><w:r><w:t>About Me:</w:t></w:r></w:p>"
$body += "<w:p><w:r><w:tab/><w:t>Background and Interests</w:t></w:r><w:r><w:tab/></w:r></w:p>"
$body += "<w:p><w:pPr><w:ind w:firstLine=`"720`"/></w:pPr><w:r><w:t>Born in New Delhi, India, I came to the U.S. in the early 90s.  My parents, siblings, and I have lived in the DMV area ever since.  I also completed all my schooling locally.</w:t></w:r></w:p>"
$body += "<w:p><w:r><w:t xml:space=`"preserve`">" + "I" + $apos + "m an avid practitioner of Yoga and Meditation and enjoy reading up on the subject during my free time.  I developed this interest after a friend pointed me towards them during a difficult time in my life and I" + $apos + "ve been hooked ever since given how much it helped me.  " + "</w:t></w:r>"+"<w:r><w:t xml:space=`"preserve`">" + "Using what I learn during the boot camp, " + "</w:t></w:r>"+"<w:r><w:t xml:space=`"preserve`">" + "I would like to create a website where I can share my development as a Yoga practitioner and share some tips or materials on Meditation. " + "</w:t></w:r></w:p>"
$body += "<w:p><w:pPr><w:ind w:firstLine=`"720`"/></w:pPr><w:r><w:t>Career:</w:t></w:r></w:p>"
$body += "<w:p><w:pPr><w:ind w:firstLine=`"720`"/></w:pPr>" + "<w:r><w:t xml:space=`"preserve`">" + "My career started " + "</w:t></w:r>" + "<w:r><w:t>" + "off" + "</w:t></w:r>" + "<w:r><w:t xml:space=`"preserve`">" + " in the mortgage industry and I have stayed in it the entire time.  " + "</w:t></w:r>" + "<w:r><w:t>" + "I" + $apos + "m a senior valuation analyst and m" + "</w:t></w:r>" + "<w:r><w:t xml:space=`"preserve`">" + "y expertise " + "</w:t></w:r>" + "<w:r><w:t>" + "include, but not limited to," + "</w:t></w:r>" + "<w:r><w:t xml:space=`"preserve`">" + " the valuation of Mortgaged Backed Securities (MBS) of various types including structured and non-structured products.  " + "</w:t></w:r>" + "<w:r><w:t xml:space=`"preserve`">" + "I" + $apos + "ve always had an analytical mindset and been interested in data " + "</w:t></w:r>" + "<w:r><w:t>" + "analytics," + "</w:t></w:r>" + "<w:r><w:t xml:space=`"preserve`">" + " so my current job suits my interests." + "</w:t></w:r>" + "</w:p>"
$body += "<w:p><w:r><w:t xml:space=`"preserve`">" + "My job required the use of various tools such as MS Office Suite, Bloomberg, " + "</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Intex</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">" + ", and Tableau to name a few.  Being in a senior and leadership role, I am also involved in various initiatives related to the future state of our organization.  I" + $apos + "m also an active member of the Diversity and Inclusion, D&amp;I, committee in our organization as I believe equal opportunity and fair treatment at work is a must for all.  " + "</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$body += "<w:p/><w:p/>"
$body += "<w:p><w:r><w:tab/><w:t xml:space=`"preserve`">   </w:t></w:r></w:p>"
$body += "<w:p><w:r><w:tab/></w:r></w:p>"

$d = $word.ActiveDocument

$d.Content.InsertXML($body)

$listStyle = $d.Styles("ListParagraph")
if ($listStyle) { $listStyle.Delete() }

